$d = $word.ActiveDocument
$t = $d.Tables(1)
$values = @(
    "63-13=50",
    "98-27=71",
    "21+34=55",
    "68-43=25",
    "79-69=10",
    "86+3=89",
    "84+10=94",
    "71-26=45",
    "13+84=97",
    "85-62=23",
    "20-17=3",
    "65+4=69",
    "86-76=10",
    "9+82=91",
    "67-12=55",
    "57-18=39",
    "69-48=21",
    "32+53=85",
    "86-10=76",
    "38-4=34",
    "37-7=30",
    "83-61=22",
    "67-14=53",
    "16+83=99",
    "69-66=3",
    "77-16=61",
    "27+37=64",
    "71-65=6",
    "43-4=39",
    "60-17=43",
    "53-36=17",
    "19+39=58",
    "9+52=61",
    "84-44=40",
    "2+35=37",
    "85-84=1",
    "27+61=88",
    "8+33=41",
    "65+2=67",
    "14+55=69",
    "51-26=25",
    "51-35=16",
    "57+26=83",
    "73-34=39",
    "26+61=87",
    "81-60=21",
    "16+57=73",
    "59+35=94",
    "47-16=31",
    "59+26=85",
    "51-45=6",
    "23-2=21",
    "35+61=96",
    "91-57=34",
    "47-25=22",
    "26-3=23",
    "99-45=54",
    "3+64=67",
    "79-28=51",
    "71-54=17",
    "42-19=23",
    "12+63=75",
    "83+16=99",
    "16+63=79",
    "95-91=4",
    "0+14=14",
    "2+55=57",
    "21+23=44",
    "24+72=96",
    "68+23=91",
    "65-38=27",
    "78-35=43",
    "49+0=49",
    "32+10=42",
    "35-28=7",
    "31+29=60",
    "4+63=67",
    "15+14=29",
    "98-94=4",
    "48-6=42",
    "22+1=23",
    "6+45=51",
    "76-23=53",
    "97-86=11",
    "68-64=4",
    "86-83=3",
    "24+50=74",
    "52+26=78",
    "98-15=83",
    "44+7=51",
    "89-78=11",
    "54-50=4",
    "33-13=20",
    "36+57=93",
    "48+46=94",
    "99-46=53",
    "96-90=6",
    "1+35=36",
    "91-48=43",
    "51+33=84"
)
$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}
Write-Host "Done. Updated" $idx "cells."
